$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update Price (D) and Volume(1h) (E) for rows 2-34 ---
$ws.Range("D2").Value = "25.775.67"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.745.34"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'235.25"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5077"
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("D8").Value = "'40.42"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").Value = "'0.2663"
$ws.Range("E9").Value = "  +4.78%  "
$ws.Range("D10").Value = "'0.06163"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("D11").Value = "1.738.29"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'0.06938"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "'15.22"
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").Value = "'0.6206"
$ws.Range("E14").Value = "  +9.63%  "
$ws.Range("D15").Value = "'4.463"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "'77.73"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "25.813.00"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "'11.56"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("D21").Value = "'0.000006624"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "1.962.60"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'4.034"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'8.231"
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("D25").Value = "'5.119"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("D26").Value = "'136.43"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").Value = "'1.461"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").Value = "'14.99"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").Value = "'1.755"
$ws.Range("E29").Value = "  -3.49%  "
$ws.Range("D30").Value = "'102.33"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'0.08121"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "'3.674"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").Value = "'3.378"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'0.04384"
$ws.Range("E34").Value = "  +0.06%  "

# --- Step 2: insert new row for Frax at row 35, shifting old rows 35-51 down to 36-52 ---
$ws.Range("A35:E35").Insert()

# Carry the "index" column formatting (border/bold/centered) down onto the new row,
# matching the look of every other data row in column A.
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: populate the newly inserted Frax row ---
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = "  +0.00%  "

# --- Step 4: fix the A-column sequence numbers for shifted rows 36-51 ---
$ws.Range("A36").Value = 34
$ws.Range("A37").Value = 35
$ws.Range("A38").Value = 36
$ws.Range("A39").Value = 37
$ws.Range("A40").Value = 38
$ws.Range("A41").Value = 39
$ws.Range("A42").Value = 40
$ws.Range("A43").Value = 41
$ws.Range("A44").Value = 42
$ws.Range("A45").Value = 43
$ws.Range("A46").Value = 44
$ws.Range("A47").Value = 45
$ws.Range("A48").Value = 46
$ws.Range("A49").Value = 47
$ws.Range("A50").Value = 48
$ws.Range("A51").Value = 49

# --- Step 5: update Price/Volume for the shifted rows (old rows 35-50, now at 36-51) ---
$ws.Range("D36").Value = "'2.647"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("D37").Value = "'0.9921"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "'0.5988"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "'2.603"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").Value = "'0.01553"
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("D41").Value = "'1.918"
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'101.59"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "'0.3809"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").Value = "'0.7465"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").Value = "'4.887"
$ws.Range("E46").Value = "  -5.26%  "
$ws.Range("D47").Value = "'0.05493"
$ws.Range("E47").Value = "  +4.89%  "
$ws.Range("D48").Value = "'0.1094"
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("D49").Value = "'5.897"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "'29.96"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "'52.37"
$ws.Range("E51").Value = "  +0.24%  "

# --- Step 6: remove the now-duplicated last row (old USDD, shifted to row 52) ---
$ws.Range("A52:E52").Delete()
